$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: stash copies of every distinct cell style we will need to re-apply
# later into out-of-the-way "scratch" cells, BEFORE any target cell is
# touched. Several donor cells are themselves later overwritten with new
# content/formatting, so their original style must be preserved first.
# ---------------------------------------------------------------------------

# style "s1"  (quotePrefix, default font) - donor: I5
$ws.Range("I5").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)

# style "s3"  (quotePrefix, applyFont) - donor: H4
$ws.Range("H4").Copy()
$ws.Range("A21").PasteSpecial($xlPasteFormats)

# style "s4"  (applyFont, default font) - donor: B4 (stays s4 unchanged, safe)
$ws.Range("B4").Copy()
$ws.Range("A22").PasteSpecial($xlPasteFormats)

# style "s5"  (underlined font) - donor: A5
$ws.Range("A5").Copy()
$ws.Range("A23").PasteSpecial($xlPasteFormats)

# style "s6"  (underlined font, quotePrefix) - donor: H6
$ws.Range("H6").Copy()
$ws.Range("A24").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: apply the stashed styles to every target cell that needs them,
# using the scratch cells (never themselves overwritten) as stable donors.
# Cells whose style does not change (A6, H5, H7, I7) are intentionally left
# untouched.
# ---------------------------------------------------------------------------

# s1 -> H4, I4, H5 (unchanged, skip), H7/I7 (unchanged, skip), H9, I4, I8, I9, H10
$ws.Range("A20").Copy()
$ws.Range("H4").PasteSpecial($xlPasteFormats)
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("I8").PasteSpecial($xlPasteFormats)
$ws.Range("H9").PasteSpecial($xlPasteFormats)
$ws.Range("I9").PasteSpecial($xlPasteFormats)
$ws.Range("H10").PasteSpecial($xlPasteFormats)

# s3 -> H6
$ws.Range("A21").Copy()
$ws.Range("H6").PasteSpecial($xlPasteFormats)

# s4 -> A8, A9, A10  (A6 already s4 and unchanged; A4 loses s4, see step 3)
$ws.Range("A22").Copy()
$ws.Range("A8").PasteSpecial($xlPasteFormats)
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("A10").PasteSpecial($xlPasteFormats)

# s5 -> A7, C10, G14
$ws.Range("A23").Copy()
$ws.Range("A7").PasteSpecial($xlPasteFormats)
$ws.Range("C10").PasteSpecial($xlPasteFormats)
$ws.Range("G14").PasteSpecial($xlPasteFormats)

# s6 -> H8
$ws.Range("A24").Copy()
$ws.Range("H8").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# remove the scratch helper cells entirely (value + format) so they don't
# leak into the saved worksheet / expand its used range.
$ws.Range("A20:A24").Clear()

# ---------------------------------------------------------------------------
# Step 3: cells that must end up with NO style attribute (plain/default).
# ---------------------------------------------------------------------------
$ws.Range("A4").ClearFormats()
$ws.Range("A5").ClearFormats()

# cells that disappear entirely from the target layout.
$ws.Range("I5").Clear()
$ws.Range("I6").Clear()
$ws.Range("G10").Clear()

# ---------------------------------------------------------------------------
# Step 4: write the new/changed cell values. Numeric-looking strings are
# prefixed with a single quote so Excel stores them as text (shared string)
# instead of converting them to numbers, matching the source data.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "REG-202"
$ws.Range("H4").Value = "'232435343"
$ws.Range("I4").Value = "'6"

$ws.Range("A5").Value = "REG-203"
$ws.Range("H5").Value = "'2232444"

$ws.Range("A6").Value = "REG-401"
$ws.Range("H6").Value = "'12000000"

$ws.Range("A7").Value = "REG-402"
$ws.Range("H7").Value = "'3250000"
$ws.Range("I7").Value = "'5"

$ws.Range("A8").Value = "REG-403"
$ws.Range("H8").Value = "'6963455"
$ws.Range("I8").Value = "'2"

$ws.Range("A9").Value = "REG-404"
$ws.Range("C9").Value = "NullPointerException"
$ws.Range("H9").Value = "'543453"
$ws.Range("I9").Value = "'4"

$ws.Range("A10").Value = "REG-501"
$ws.Range("H10").Value = "'23243543"

# C10 and G14 stay blank (format-only cells) -- already styled in Step 2.

# ---------------------------------------------------------------------------
# Step 5: selection / active cell, matching the final authored state.
# ---------------------------------------------------------------------------
$ws.Range("D13").Select()
